$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.412.17'
$ws.Range('E2').Value = '  +3.41%  '
$ws.Range('D3').Value = '1.868.56'
$ws.Range('E3').Value = '  +1.96%  '
$ws.Range('D5').Value = "'337.93"
$ws.Range('E5').Value = '  +1.91%  '
$ws.Range('D6').Value = "'1.002"
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').Value = "'0.4687"
$ws.Range('E7').Value = '  +1.56%  '
$ws.Range('D8').Value = "'0.3966"
$ws.Range('E8').Value = '  +3.38%  '
$ws.Range('D9').Value = "'47.62"
$ws.Range('E9').Value = '  +2.11%  '
$ws.Range('D10').Value = "'0.08027"
$ws.Range('E10').Value = '  +1.63%  '
$ws.Range('D11').Value = "'0.9991"
$ws.Range('E11').Value = '  +2.78%  '
$ws.Range('D12').Value = "'21.98"
$ws.Range('D13').Value = "'6.040"
$ws.Range('E13').Value = '  +2.57%  '
$ws.Range('D14').Value = '1.866.65'
$ws.Range('E14').Value = '  +1.25%  '
$ws.Range('D15').Value = "'7.261"
$ws.Range('E15').Value = '  +2.86%  '
$ws.Range('D16').Value = "'90.75"
$ws.Range('E16').Value = '  +3.04%  '
$ws.Range('D17').Value = "'1.003"
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('E18').Value = '  +1.26%  '
$ws.Range('D19').Value = "'0.06628"
$ws.Range('D20').Value = "'17.51"
$ws.Range('E20').Value = '  +2.76%  '
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('D22').Value = '28.413.48'
$ws.Range('E22').Value = '  +3.40%  '
$ws.Range('D23').Value = "'5.471"
$ws.Range('E23').Value = '  +2.42%  '
$ws.Range('D24').Value = "'11.05"
$ws.Range('E24').Value = '  +2.18%  '
$ws.Range('D25').Value = "'2.267"
$ws.Range('E25').Value = '  -1.88%  '
$ws.Range('D26').Value = '2.086.43'
$ws.Range('E26').Value = '  +1.34%  '
$ws.Range('D27').Value = "'160.62"
$ws.Range('E27').Value = '  +2.06%  '
$ws.Range('D28').Value = "'19.74"
$ws.Range('E28').Value = '  +1.62%  '
$ws.Range('D29').Value = "'2.116"
$ws.Range('E29').Value = '  +2.43%  '
$ws.Range('D30').Value = "'5.492"
$ws.Range('E30').Value = '  +3.61%  '
$ws.Range('D31').Value = "'120.11"
$ws.Range('E31').Value = '  +0.92%  '
$ws.Range('D32').Value = "'0.9702"
$ws.Range('E32').Value = '  +1.35%  '
$ws.Range('D33').Value = "'0.09491"
$ws.Range('E33').Value = '  +2.05%  '
$ws.Range('D34').Value = "'3.590"
$ws.Range('E34').Value = '  +0.44%  '
$ws.Range('D35').Value = "'5.356"
$ws.Range('E35').Value = '  +2.25%  '
$ws.Range('D36').Value = "'1.374"
$ws.Range('E36').Value = '  +4.36%  '
$ws.Range('D37').Value = "'0.06101"
$ws.Range('E37').Value = '  +2.64%  '
$ws.Range('E38').Value = '  +2.13%  '
$ws.Range('D39').Value = "'8.330"
$ws.Range('E39').Value = '  +3.16%  '
$ws.Range('D40').Value = "'1.179"
$ws.Range('E40').Value = '  +2.04%  '
$ws.Range('D41').Value = "'0.5937"
$ws.Range('E41').Value = '  +2.26%  '
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('E43').Value = '  +1.77%  '
$ws.Range('E44').Value = '  +3.17%  '
$ws.Range('E45').Value = '  -0.34%  '
$ws.Range('D46').Value = "'0.5564"
$ws.Range('E46').Value = '  +1.51%  '
$ws.Range('D47').Value = "'12.13"
$ws.Range('E47').Value = '  +1.23%  '
$ws.Range('D48').Value = "'1.954"
$ws.Range('E48').Value = '  +4.30%  '
$ws.Range('D49').Value = "'0.06989"
$ws.Range('E49').Value = '  +5.21%  '
$ws.Range('D50').Value = "'2.065"
$ws.Range('E50').Value = '  +13.91%  '
$ws.Range('D51').Value = "'111.65"
$ws.Range('E51').Value = '  +1.09%  '
